# Signup and quiz status
# - Populate the "harsha" worksheet (3rd tab) with signup/quiz data
# - Add mailto: hyperlinks on the email column
# - Make "harsha" the active/selected tab (instead of "manoj")

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(3)

# Header row
$ws.Range("A1").Value = "name"
$ws.Range("B1").Value = "username"
$ws.Range("C1").Value = "email"
$ws.Range("D1").Value = "password"
$ws.Range("E1").Value = "role"
$ws.Range("F1").Value = "skills"

# Row 2 - harsha (HR)
$ws.Range("A2").Value = "harsha1222222"
$ws.Range("B2").Value = "harsha12344690"
$ws.Range("C2").Value = "harsha12389@yopmail.com"
$ws.Range("D2").Value = """12345678"""
$ws.Range("E2").Value = "HR"

# Row 3 - vardhan (Subject Expert)
$ws.Range("A3").Value = "vardhan123444"
$ws.Range("B3").Value = "vardhan12356666"
$ws.Range("C3").Value = "vardhan124444@yopmail.com"
$ws.Range("D3").Value = """12345678"""
$ws.Range("E3").Value = "Subject Expert"
$ws.Range("F3").Value = "html"
$ws.Range("G3").Value = "css"

# Hyperlinks on the email cells
$ws.Hyperlinks.Add($ws.Range("C2"), "mailto:harsha12389@yopmail.com", "", "", "harsha12389@yopmail.com") | Out-Null
$ws.Hyperlinks.Add($ws.Range("C3"), "mailto:vardhan124444@yopmail.com", "", "", "vardhan124444@yopmail.com") | Out-Null

# Column widths (approximate autofit widths; COM ColumnWidth only has 1/7-character
# granularity so these are the closest reachable values to the saved widths of
# 15.54296875 / 16.26953125 / 24.08984375 / 25 / 15.26953125)
$ws.Columns.Item(1).ColumnWidth = 14.714285714285714
$ws.Columns.Item(2).ColumnWidth = 15.428571428571429
$ws.Columns.Item(3).ColumnWidth = 23.285714285714285
$ws.Columns.Item(4).ColumnWidth = 24.142857142857142
$ws.Columns.Item(5).ColumnWidth = 14.428571428571429

# Select C2 as the active cell on this sheet
$ws.Range("C2").Select() | Out-Null

# Make "harsha" the active tab
$ws.Activate() | Out-Null

Write-Host "done"
